$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking
# values (e.g. "0.999", "8.42") are stored as text, not converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.113.16"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "2.563.33"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D5").Value = "583.88"
$ws.Range("E5").Value = "  +2.58%  "
$ws.Range("D6").Value = "147.45"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("E9").Value = "  +2.79%  "
$ws.Range("D10").Value = "5.63"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").Value = "27.37"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "3.024.60"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").Value = "63.059.83"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").Value = "2.530.45"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "11.34"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("D19").Value = "343.42"
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("D20").Value = "4.42"
$ws.Range("E20").Value = "  +2.69%  "
$ws.Range("D21").Value = "6.87"
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  -4.02%  "
$ws.Range("D24").Value = "66.72"
$ws.Range("D25").Value = "2.695.63"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "0.170"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "1.62"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").Value = "8.09"
$ws.Range("E28").Value = "  +9.55%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "8.42"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("E32").Value = "  +7.17%  "
$ws.Range("D33").Value = "0.0₃0823"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").Value = "461.40"
$ws.Range("E34").Value = "  +12.23%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "175.61"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "1.62"
$ws.Range("E36").Value = "  +2.51%  "
$ws.Range("D37").Value = "0.408"
$ws.Range("E37").Value = "  +2.23%  "
$ws.Range("D38").Value = "19.18"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").Value = "4.53"
$ws.Range("E39").Value = "  +2.87%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("D43").Value = "150.91"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").Value = "3.82"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("D45").Value = "20.85"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("E46").Value = "  +4.78%  "
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("D48").Value = "0.0974"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("E50").Value = "  -2.78%  "
$ws.Range("E51").Value = "  +0.28%  "
